# modify the childnodes counting and simplify the codes
#
# The GO-term "id" column got overwritten with the literal text of the
# unbound `id` builtin (a bug from simplifying the Python scraper code:
# `row.id` became just `id`, i.e. `str(id)` -> "<built-in function id>"),
# and the childnodes counts (column D) were recomputed with the new
# counting logic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$idText = "<built-in function id>"

# New childnodes counts (column D) after the counting-logic change.
$childCounts = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 3
    6  = 0
    7  = 2
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 1
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 0
    24 = 1
    25 = 0
    26 = 0
    27 = 0
    28 = 2
    29 = 0
    30 = 0
    31 = 0
    32 = 0
    33 = 3
    34 = 0
    35 = 0
    36 = 2
}

for ($row = 2; $row -le 36; $row++) {
    $ws.Cells.Item($row, 1).Value = $idText
    $ws.Cells.Item($row, 4).Value = $childCounts[$row]
}
